# Added loop for prediction draws - update the aggregate calibration output
# values for row 2 (the single simulation draw row) with the new results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10801999.03673311
$ws.Range("C2").Value = 2414466.439874003
$ws.Range("D2").Value = 28111992.96818429
$ws.Range("E2").Value = 1189901.79008772
$ws.Range("F2").Value = 9486554.436699318
$ws.Range("G2").Value = 1806289.812291691
$ws.Range("H2").Value = 2157083.447618817
$ws.Range("I2").Value = 10801999.03673311
$ws.Range("J2").Value = 45425814
$ws.Range("L2").Value = 30526459.40805829
$ws.Range("M2").Value = 10676456.22678704
$ws.Range("N2").Value = 3963373.259910508
$ws.Range("O2").Value = 38154.79029512414
$ws.Range("P2").Value = 198626.9279718958
$ws.Range("Q2").Value = 236781.71826702
$ws.Range("R2").Value = 361.6488484323222
$ws.Range("S2").Value = 112369.8360000918
$ws.Range("T2").Value = 112731.4848485241
